$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 (title slide): update the title text and let the autofit shrink
# ("fontScale") recompute/clear since the new text fits without scaling.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Lecture 14:LLM Prompting"
# Force a normAutofit recompute (drops the stale fontScale="90000").
$titleShape.TextFrame.AutoSize = 2

# ---------------------------------------------------------------------------
# Slide 2: rewrite the bullet list content (Vision Transformer topics ->
# LLM Prompting topics).
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item(2)

# Paragraph text (in order) and desired indent level (1 = top level,
# 2 = first sub-level, matching PowerPoint's 1-based IndentLevel).
$paraTexts  = @(
    "Prompting Large Language Model",
    "CoT",
    "Self-Consistency",
    "ToT",
    "React",
    "PoT",
    "Least to Most",
    "Self-Refine",
    "Self-Ask",
    "",
    ""
)
$paraLevels = @(1, 2, 2, 2, 2, 2, 2, 2, 2, 1, 1)

$fullText = [string]::Join("`r", $paraTexts)
$contentShape.TextFrame.TextRange.Text = $fullText
$contentShape.TextFrame.AutoSize = 2

$tr = $contentShape.TextFrame.TextRange
$pos = 1
for ($i = 0; $i -lt $paraTexts.Length; $i++) {
    $len = $paraTexts[$i].Length
    if ($len -gt 0 -and $paraLevels[$i] -gt 1) {
        $run = $tr.Characters($pos, $len)
        $run.IndentLevel = $paraLevels[$i]
    }
    $pos += $len + 1
}
